$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H) added to the s_vals sheet.
# Header H1 should carry the same formatting (bold, border, centered) as
# the other header cells, so copy formats from G1 (xlPasteFormats = -4122)
# before writing the header text.
$xlPasteFormats = -4122
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial($xlPasteFormats)
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for rows 2-5 (plain numbers, default style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
